$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "...Soundscape Narrative, Visual Argument, Website, and Consolidation;"
#    -> "...Soundscape Narrative (mp3), Visual Argument (png), Website
#        (html), and Consolidation (tbd);"
#    Insert each "(ext)" right after its project name so we touch only the
#    run that already contains that word, instead of re-writing the whole
#    sentence (which would swallow neighbouring runs on save).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Soundscape Narrative", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" (mp3)")

$rng = $d.Content
$rng.Find.Execute("Visual Argument", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" (png)")

$rng = $d.Content
$rng.Find.Execute("Website", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" (html)")

$rng = $d.Content
$rng.Find.Execute("Consolidation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" (tbd)")

# ---------------------------------------------------------------------
# 2. ", hyperlinked to a point " -> ", hyperlinked to a commit "
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "hyperlinked to a point", $true, $false, $false, $false, $false, $true, 1, $false,
    "hyperlinked to a commit", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Drop the line break inside the NB parenthetical and shrink its font
#    to 8pt (sz/szCs = 16 half-points).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("instead. ^l" + "But I hope", $true, $false, $false, $false, $false, $true, 1, $false,
    "instead. But I hope", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute(
    "(NB: if you prefer your refection stay private, you may email me instead. But I hope you" + [char]8217 + "ll consider sharing your reflections with your peers.)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.Size = 8

# ---------------------------------------------------------------------
# 4. "...words in total, but your mileage..." -> "...words in total (800
#    at a minimum), but your mileage..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("in total", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" (800 at a minimum)")

# ---------------------------------------------------------------------
# 5. Replace the due-date line in the "first page" header of section 1
#    (header3.xml / wdHeaderFooterFirstPage).
# ---------------------------------------------------------------------
$hdr = $d.Sections(1).Headers(2)
$hdr.Range.Find.Execute(
    "Due Tuesday, April 21, at 1:50pm",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Due in our final exam slot (see course site)",
    2) | Out-Null
